$d = $word.ActiveDocument
$bmBody = $d.Bookmarks("BODY")
$rng = $d.Range($bmBody.Start, $bmBody.End)
$d.Bookmarks.Add("_GoBack", $rng)
# now re-add BODY fresh at position 0 pointing to same full range, maybe it reorders
$bmBody2 = $d.Bookmarks("BODY")
$rng2 = $d.Range($bmBody2.Start, $bmBody2.End)
$d.Bookmarks.Add("BODY", $rng2)
